$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.882.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '''2.622.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''594.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = '''151.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.588'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  +4.74%  '
$ws.Range("D10").Value = '''5.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  +2.83%  '
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").Value = '''27.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").Value = '''3.091.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000170'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +13.96%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '''63.744.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '''2.605.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = '''12.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = '''4.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = '''348.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '''7.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.00%  '
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").Value = '''67.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("D24").Value = '''1.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.94%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''9.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '''1.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("D27").Value = '''8.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").Value = '''544.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.60%  '
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '''0.0₃0904'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.09%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("E33").Value = '  +4.96%  '
$ws.Range("D34").Value = '''5.46'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.28%  '
$ws.Range("D35").Value = '''6.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").Value = '''0.424'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.90%  '
$ws.Range("D37").Value = '''166.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("D38").Value = '''20.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '''169.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("E43").Value = '  +4.40%  '
$ws.Range("D44").Value = '''23.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.21%  '
$ws.Range("E45").Value = '  -2.40%  '
$ws.Range("D46").Value = '''2.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.38%  '
$ws.Range("D47").Value = '''0.637'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").Value = '''0.0253'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").Value = '''0.0972'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("D50").Value = '''19.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").Value = '''0.0₆0233'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +20.23%  '
